# InvoiceLogTemplate.xlsx - "Fix Code Bugs and Add Field For Invoice Generator"
#
# Content-level changes applied:
#   1. On the "Client_List" sheet, fill column G (rows 2-14) with the new
#      "demo adress" field value (new shared string, referenced 13 times -
#      this is the new field added for the invoice generator).
#   2. Move the active/selected sheet + cell from "Project_List"!F10 to
#      "Client_List"!K10 (the workbook was left with Client_List as the
#      active tab after the edits were made).

$wb = $excel.ActiveWorkbook

# 1. Add the new "demo adress" field to Client_List, column G, rows 2-14.
$wsClientList = $wb.Worksheets.Item("Client_List")
$wsClientList.Range("G2:G14").Value = "demo adress"

# 2. Re-point the active sheet/selection: Client_List becomes the active tab
#    with K10 selected; Project_List keeps its own F10 selection but is no
#    longer the tab shown on open.
$wsClientList.Activate()
$wsClientList.Range("K10").Select()
